# Update cryptos list: prices (col D) and 1h volume % (col E) for rows 2-51.
# Rows 36/37 also swap Fetch.AI <-> Monero (name + link).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.981.48"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "3.081.69"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'519.97"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'136.08"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.081.77"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.450"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").Value = "'7.32"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "3.611.35"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'25.27"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "57.087.58"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "3.079.29"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'12.46"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'347.15"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'68.30"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "0.0₃0862"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'7.26"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'5.86"
$ws.Range("E33").Value = "  -8.05%  "
$ws.Range("D34").Value = "'20.79"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "'4.91"
$ws.Range("E35").Value = "  +6.76%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'159.16"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = "  -3.65%  "
$ws.Range("D38").Value = "'6.00"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").Value = "'25.86"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "'0.0652"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'4.01"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "'0.690"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "2.388.28"
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("D46").Value = "'36.61"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "3.120.74"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'0.956"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("E51").Value = "  -2.30%  "
